{"js": "// Fix systematic spacing issue between header bar and body text\n// - Collapse the three long CORE COMPETENCIES detail paragraphs into a single\n//   summary line.\n// - Add a new \"TECHNICAL SKILLS\" section (heading + three detail paragraphs)\n//   at the end of the document.\n\nconst body = context.document.body;\nbody.load(\"paragraphs/items/text\");\nawait context.sync();\n\nconst paras = body.paragraphs.items;\n\n// Locate the three CORE COMPETENCIES detail paragraphs by their distinctive\n// leading text rather than a hard-coded index, so the script is resilient to\n// minor structural differences.\nlet coreIdx = -1;\nfor (let i = 0; i < paras.length; i++) {\n  if (paras[i].text.indexOf(\"Data Visualization & Design: Interactive Dashboards\") === 0) {\n    coreIdx = i;\n    break;\n  }\n}\n\nif (coreIdx !== -1) {\n  // Replace the first of the three paragraphs with the condensed summary...\n  paras[coreIdx].insertText(\n    \"Data Visualization & Design \\u2022 Geospatial Analysis & Mapping \\u2022 Technical Visualization\",\n    \"Replace\"\n  );\n  // ...and delete the next two (the Geospatial Analysis & Mapping and\n  // Technical Visualization detail paragraphs).\n  paras[coreIdx + 1].delete();\n  paras[coreIdx + 2].delete();\n}\n\n// Append the new TECHNICAL SKILLS section at the end of the document body.\n// Note: create all four paragraphs first, THEN set the heading's style \u2014\n// setting the style immediately after creating the heading paragraph would\n// cause the following (End-anchored) paragraphs to inherit \"Heading2\" since\n// insertParagraph(\"End\") picks up formatting from the current last paragraph.\nconst heading = body.insertParagraph(\"TECHNICAL SKILLS\", \"End\");\n\nbody.insertParagraph(\n  \"DATA VISUALIZATION & DESIGN Interactive Dashboards; Statistical Visualization; Geospatial Mapping; Choropleth Design; Web Visualization; Presentation Design; Data Storytelling\",\n  \"End\"\n);\n\nbody.insertParagraph(\n  \"GEOSPATIAL ANALYSIS & MAPPING Spatial Analysis; Mapping Technologies; Web Mapping; Spatial Data Processing; Census Data Integration; Custom Tile Servers; Spatial Clustering\",\n  \"End\"\n);\n\nbody.insertParagraph(\n  \"TECHNICAL VISUALIZATION Programming; Database Integration; Cloud Platforms; Web Technologies; Statistical Computing; Version Control; DevOps\",\n  \"End\"\n);\n\nheading.style = \"Heading 2\";\n\nawait context.sync();\n", "ps1": "# Fix systematic spacing issue between header bar and body text\n# - Collapse the three long CORE COMPETENCIES detail paragraphs into a single\n#   summary line.\n# - Add a new \"TECHNICAL SKILLS\" section (heading + three detail paragraphs)\n#   at the end of the document.\n\n$d = $word.ActiveDocument\n\n# Locate the first of the three CORE COMPETENCIES detail paragraphs by its\n# distinctive leading text rather than a hard-coded index, so the script is\n# resilient to minor structural differences.\n$coreIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text\n  if ($t.StartsWith(\"Data Visualization & Design: Interactive Dashboards\")) {\n    $coreIdx = $i\n    break\n  }\n}\n\nif ($coreIdx -ne -1) {\n  # Replace the first of the three paragraphs with the condensed summary...\n  $d.Paragraphs.Item($coreIdx).Range.Text = \"Data Visualization & Design \u2022 Geospatial Analysis & Mapping \u2022 Technical Visualization\"\n  # ...and delete the next two (the Geospatial Analysis & Mapping and\n  # Technical Visualization detail paragraphs). Deleting the paragraph at\n  # $coreIdx + 1 twice removes both, since each delete shifts later\n  # paragraphs up by one.\n  $d.Paragraphs.Item($coreIdx + 1).Range.Delete()\n  $d.Paragraphs.Item($coreIdx + 1).Range.Delete()\n}\n\n# Append the new TECHNICAL SKILLS section at the end of the document body.\n# Note: insert all four paragraphs first, THEN set the heading's style --\n# setting the style immediately after creating the heading paragraph would\n# cause the following paragraphs (each inserted after the new last\n# paragraph) to inherit the \"Heading2\" style.\n$lastIdx = $d.Paragraphs.Count\n$d.Paragraphs.Item($lastIdx).Range.InsertParagraphAfter()\n\n$headingIdx = $d.Paragraphs.Count\n$d.Paragraphs.Item($headingIdx).Range.Text = \"TECHNICAL SKILLS\"\n\n$d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertParagraphAfter()\n$d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = \"DATA VISUALIZATION & DESIGN Interactive Dashboards; Statistical Visualization; Geospatial Mapping; Choropleth Design; Web Visualization; Presentation Design; Data Storytelling\"\n\n$d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertParagraphAfter()\n$d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = \"GEOSPATIAL ANALYSIS & MAPPING Spatial Analysis; Mapping Technologies; Web Mapping; Spatial Data Processing; Census Data Integration; Custom Tile Servers; Spatial Clustering\"\n\n$d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertParagraphAfter()\n$d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = \"TECHNICAL VISUALIZATION Programming; Database Integration; Cloud Platforms; Web Technologies; Statistical Computing; Version Control; DevOps\"\n\n$d.Paragraphs.Item($headingIdx).Style = \"Heading 2\"\n"}
